$d = $word.ActiveDocument

# 1. "...statistical methods that fully utilize the..." ->
#    "...statistical methods that leverage the..."
$d.Content.Find.Execute("fully utilize the", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "leverage the", 2)

# 2. "...quantitative merits of either have not..." ->
#    "...quantitative merits have not..."
$d.Content.Find.Execute("quantitative merits of either have not", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "quantitative merits have not", 2)

# 3. "...Patuxent River estuary..." -> "...Patuxent River Estuary..."
$d.Content.Find.Execute("Patuxent River estuary", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Patuxent River Estuary", 2)

# 4. Reposition the "_GoBack" bookmark: in the signature block it currently
#    sits right after "Marcus W. Beck" (end of document). Move it so it
#    sits between the 2nd and 3rd tab characters that precede the
#    signature name, matching the author's last edit location.
$sigPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$sigStart = $sigPara.Range.Start

try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
} catch {}

$newPos = $sigStart + 2
$newRange = $d.Range($newPos, $newPos)
$d.Bookmarks.Add("_GoBack", $newRange)
